# Update a handful of numeric values in column F across the four sheets
# of the "上海-漫展信息" workbook, matching the regenerated data pull.

$wb = $excel.ActiveWorkbook

# Sheet 1: "展览"
$ws1 = $wb.Worksheets.Item(1)
$ws1.Cells.Item(5, 6).Value  = 8162
$ws1.Cells.Item(14, 6).Value = 1171
$ws1.Cells.Item(16, 6).Value = 22
$ws1.Cells.Item(17, 6).Value = 746
$ws1.Cells.Item(23, 6).Value = 6858
$ws1.Cells.Item(25, 6).Value = 53982
$ws1.Cells.Item(26, 6).Value = 4199
$ws1.Cells.Item(27, 6).Value = 2
$ws1.Cells.Item(29, 6).Value = 802
$ws1.Cells.Item(35, 6).Value = 2037
$ws1.Cells.Item(38, 6).Value = 838
$ws1.Cells.Item(39, 6).Value = 1080
$ws1.Cells.Item(40, 6).Value = 472
$ws1.Cells.Item(47, 6).Value = 117

# Sheet 2: "演出"
$ws2 = $wb.Worksheets.Item(2)
$ws2.Cells.Item(28, 6).Value = 113
$ws2.Cells.Item(31, 6).Value = 16

# Sheet 3: "本地生活"
$ws3 = $wb.Worksheets.Item(3)
$ws3.Cells.Item(8, 6).Value  = 2332
$ws3.Cells.Item(9, 6).Value  = 9287
$ws3.Cells.Item(10, 6).Value = 1578

# Sheet 4: "全部类型"
$ws4 = $wb.Worksheets.Item(4)
$ws4.Cells.Item(7, 6).Value  = 1578
$ws4.Cells.Item(14, 6).Value = 1171
$ws4.Cells.Item(16, 6).Value = 22
$ws4.Cells.Item(17, 6).Value = 746
$ws4.Cells.Item(20, 6).Value = 6858
$ws4.Cells.Item(22, 6).Value = 53982
$ws4.Cells.Item(27, 6).Value = 4199
$ws4.Cells.Item(38, 6).Value = 838
$ws4.Cells.Item(46, 6).Value = 117
